$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match data between row 16 and row 17 ---
# (keep A, B, C, D, E, G, I, K, O, S as-is; swap F,H,J,L,M,N,P,Q,R,T,U,V)
$cols1617 = @("F","H","J","L","M","N","P","Q","R","T","U","V")
foreach ($col in $cols1617) {
    $tmp = $ws.Range("$col`16").Value2
    $ws.Range("$col`16").Value2 = $ws.Range("$col`17").Value2
    $ws.Range("$col`17").Value2 = $tmp
}

# --- Swap match data between row 22 and row 23 ---
# (keep A, B, C, D, E, K, O, S as-is; swap F,G,H,I,J,L,M,N,P,Q,R,T,U,V)
$cols2223 = @("F","G","H","I","J","L","M","N","P","Q","R","T","U","V")
foreach ($col in $cols2223) {
    $tmp = $ws.Range("$col`22").Value2
    $ws.Range("$col`22").Value2 = $ws.Range("$col`23").Value2
    $ws.Range("$col`23").Value2 = $tmp
}

# --- Append new row 42 with new match data ---
# Clone formatting (cell styles) from the previous data row so the new
# row matches the workbook's existing look (bold/centered index column,
# date-formatted date column, etc.)
$ws.Range("A41:V41").Copy()
$ws.Range("A42:V42").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A42").Value2 = 41
$ws.Range("B42").Value2 = "lebanon"
$ws.Range("C42").Value2 = "premier-league"
$ws.Range("D42").Value2 = "2023-2024"
$ws.Range("E42").Value2 = 45235.58333333334
$ws.Range("F42").Value2 = "Al Ansar"
$ws.Range("G42").Value2 = 0
$ws.Range("H42").Value2 = "Bourj FC"
$ws.Range("I42").Value2 = 2
$ws.Range("J42").Value2 = 1.57
$ws.Range("K42").Value2 = "04/11/2023 02:13"
$ws.Range("L42").Value2 = 1.88
$ws.Range("M42").Value2 = "05/11/2023 13:12"
$ws.Range("N42").Value2 = 3.59
$ws.Range("O42").Value2 = "04/11/2023 02:13"
$ws.Range("P42").Value2 = 3.42
$ws.Range("Q42").Value2 = "05/11/2023 13:12"
$ws.Range("R42").Value2 = 4.93
$ws.Range("S42").Value2 = "04/11/2023 02:13"
$ws.Range("T42").Value2 = 4.03
$ws.Range("U42").Value2 = "05/11/2023 13:12"
$ws.Range("V42").Value2 = "https://www.betexplorer.com/football/lebanon/premier-league/al-ansar-bourj/jcV4K7Hh/"
